$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New generated-report identifiers (two tracked files replace the previous
# three: one .md file keeps "Include" handoff semantics, the other .md file
# replaces the two stale .png dependency rows from the old report).
# ---------------------------------------------------------------------------
$md1 = "8ea2fb91-5590-4825-a61c-44512c317414.md"
$md2 = "91dc6e7a-5d62-4042-b625-c13f0e55c52f.md"

$xlf1ZhCn = "8ea2fb91-5590-4825-a61c-44512c317414.71e769f1080802c231fd85751133759bd9a390b2.zh-cn.xlf"
$xlf2ZhCn = "91dc6e7a-5d62-4042-b625-c13f0e55c52f.bfeca8d639b7f0a4f621f3397be3385c8f734f64.zh-cn.xlf"
$xlf1DeDe = "8ea2fb91-5590-4825-a61c-44512c317414.71e769f1080802c231fd85751133759bd9a390b2.de-de.xlf"
$xlf2DeDe = "91dc6e7a-5d62-4042-b625-c13f0e55c52f.bfeca8d639b7f0a4f621f3397be3385c8f734f64.de-de.xlf"

$handoffTimeZhCn = "2016-03-09 18:56:18"
$handoffTimeDeDe = "2016-03-09 18:56:22"
$epoch = "0001-01-01 00:00:00"

$baseE2e    = "https://github.com/OpenLocalizationTest/oltest/blob/7fc9f7385ac446a8db4a2308410311ac6054c53d/e2e/"
$baseConfig = "https://github.com/OpenLocalizationTest/oltest/blob/7fc9f7385ac446a8db4a2308410311ac6054c53d/.localization-config"
$baseZhCnHt = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9122ca9e205039e7624f34cc15becc6d8b93a899/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$baseDeDeHt = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9e1ff68cd3b87f16ef6cab0ff58510afc7d337df/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

# ---------------------------------------------------------------------------
# Sheet "Overview": 3 columns (File Name, zh-cn, de-de); used to have a
# 4th data row (an extra .png dependency) that the new report no longer
# lists, so row 5 is dropped entirely and rows 2-3 get the new file names.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Rows.Item(5).Delete()

$wsOverview.Range("A2").Value = $md1
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

$wsOverview.Range("A3").Value = $md2
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

$wsOverview.Range("A4").Value = ".localization-config"
$wsOverview.Range("B4").Value = "Not to be localized"
$wsOverview.Range("C4").Value = "Not to be localized"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $baseE2e + $md1, "", "", $md1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $baseE2e + $md2, "", "", $md2)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $baseConfig, "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "zh-cn": 9 columns, same row-5 drop + rewrite as above, plus the
# handoff-file / handoff-datetime / status columns.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Rows.Item(5).Delete()

$wsZhCn.Range("A2").Value = $md1
$wsZhCn.Range("B2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = $xlf1ZhCn
$wsZhCn.Range("D2").Value = $handoffTimeZhCn
$wsZhCn.Range("G2").Value = $epoch
$wsZhCn.Range("H2").Value = "Include"

$wsZhCn.Range("A3").Value = $md2
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("C3").Value = $xlf2ZhCn
$wsZhCn.Range("D3").Value = $handoffTimeZhCn
$wsZhCn.Range("G3").Value = $epoch
$wsZhCn.Range("H3").Value = "Include"
$wsZhCn.Range("I3").ClearContents()

$wsZhCn.Range("A4").Value = ".localization-config"
$wsZhCn.Range("B4").Value = "Not to be localized"
$wsZhCn.Range("C4").ClearContents()
$wsZhCn.Range("D4").Value = $epoch
$wsZhCn.Range("G4").Value = $epoch
$wsZhCn.Range("H4").Value = "Ignored"
$wsZhCn.Range("I4").ClearContents()

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $baseE2e + $md1, "", "", $md1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), $baseZhCnHt + $xlf1ZhCn, "", "", $xlf1ZhCn)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $baseE2e + $md2, "", "", $md2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C3"), $baseZhCnHt + $xlf2ZhCn, "", "", $xlf2ZhCn)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), $baseConfig, "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "de-de": same shape as "zh-cn" but with the de-de handoff file
# names / timestamps / relationship targets.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Rows.Item(5).Delete()

$wsDeDe.Range("A2").Value = $md1
$wsDeDe.Range("B2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = $xlf1DeDe
$wsDeDe.Range("D2").Value = $handoffTimeDeDe
$wsDeDe.Range("G2").Value = $epoch
$wsDeDe.Range("H2").Value = "Include"

$wsDeDe.Range("A3").Value = $md2
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("C3").Value = $xlf2DeDe
$wsDeDe.Range("D3").Value = $handoffTimeDeDe
$wsDeDe.Range("G3").Value = $epoch
$wsDeDe.Range("H3").Value = "Include"
$wsDeDe.Range("I3").ClearContents()

$wsDeDe.Range("A4").Value = ".localization-config"
$wsDeDe.Range("B4").Value = "Not to be localized"
$wsDeDe.Range("C4").ClearContents()
$wsDeDe.Range("D4").Value = $epoch
$wsDeDe.Range("G4").Value = $epoch
$wsDeDe.Range("H4").Value = "Ignored"
$wsDeDe.Range("I4").ClearContents()

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $baseE2e + $md1, "", "", $md1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), $baseDeDeHt + $xlf1DeDe, "", "", $xlf1DeDe)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $baseE2e + $md2, "", "", $md2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C3"), $baseDeDeHt + $xlf2DeDe, "", "", $xlf2DeDe)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), $baseConfig, "", "", ".localization-config")
